$d = $word.ActiveDocument

# 1. Fix the typo: "being" -> "used" in the noexcept bullet.
$r = $d.Content
$r.Find.Execute("being.", $true, $false, $false, $false, $false, $true, 1, $false, "used.", 2)

# 2. Word's "_GoBack" bookmark tracks the most recent edit location. Move it from
#    the end of the previous bullet to right after the word we just typed (i.e.
#    right before the trailing period), splitting the run the same way a live
#    edit in Word would.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Force the run boundary right before "used" using a scratch bookmark...
$rBeforeUsed = $d.Content
$rBeforeUsed.Find.Execute("used", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$beforeUsed = $d.Range($rBeforeUsed.Start, $rBeforeUsed.Start)
$d.Bookmarks.Add("TempMark", $beforeUsed)

# ...then drop the real _GoBack bookmark collapsed right after "used".
$rAfterUsed = $d.Content
$rAfterUsed.Find.Execute("used", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterUsed = $d.Range($rAfterUsed.End, $rAfterUsed.End)
$d.Bookmarks.Add("_GoBack", $afterUsed)

$d.Bookmarks("TempMark").Delete()
